$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "D:\get_data_from_web"
$ws.Range("D6").Select()
